# ---------------------------------------------------------------------------
# Corrected the Sucep_new file
# Corrected the file and other functions associated with it
#
# This reshapes the susceptibility-classes grid on Sheet1:
#   - the old dilution header row (0.05, 1, 2, 4, 5, 8, 12, 16, 25, ...) is
#     replaced with text concentration labels (0.053, 0.125, 0.250, 0.500)
#     followed by the numeric series 1, 2, 4, 8, 16, 32, 64, 128
#   - the sparse 0/1/2 marker cells for each class (rows 2-14) move to new
#     columns to line up with the corrected headers
#   - the extra classes that no longer exist (rows 15-17) are removed
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Row 1 header: B1:E1 become literal text ("0.053","0.125","0.250","0.500")
#    stored as shared strings while KEEPING their existing style (s=1).
#    A plain `.Value = "0.053"` gets auto-parsed back into a number, and
#    forcing text via NumberFormat="@" (or a leading apostrophe) stamps a
#    brand-new derived style onto the cell. Routing the literal through a
#    text FORMULA first, then collapsing the formula to a static value via
#    copy / paste-special-values, yields a plain shared-string cell without
#    ever touching the style.
# ---------------------------------------------------------------------------
$headerTextCols = @(2, 3, 4, 5)
$headerTextVals = @("0.053", "0.125", "0.250", "0.500")
for ($i = 0; $i -lt $headerTextCols.Length; $i++) {
    $cell = $ws.Cells.Item(1, $headerTextCols[$i])
    $cell.Formula = "=""" + $headerTextVals[$i] + """"
    $cell.Copy()
    $cell.PasteSpecial(-4163)   # xlPasteValues
}
$excel.CutCopyMode = $false

# F1 becomes a genuine number (1) with a new "#,##0" (numFmtId 3) style.
$f1 = $ws.Cells.Item(1, 6)
$f1.NumberFormat = "#,##0"
$f1.Value = 1

# G1:J1 shift down the old dilution series (8,16 -> 2,4,8,16)
$ws.Cells.Item(1, 7).Value = 2
$ws.Cells.Item(1, 8).Value = 4
$ws.Cells.Item(1, 9).Value = 8
$ws.Cells.Item(1, 10).Value = 16
# K1:M1 (32, 64, 128) are unchanged.

# ---------------------------------------------------------------------------
# 2. Clear the marker cells that no longer apply (rows 2-14).
# ---------------------------------------------------------------------------
$cellsToClear = @(
    "K2", "L2",
    "D3", "G3", "K3",
    "J5", "L5",
    "D6", "I6", "J6",
    "C7", "M7",
    "L8", "M8",
    "K9",
    "L10",
    "G11", "H11",
    "G12", "I12", "J12",
    "I13", "K13",
    "C14"
)
foreach ($addr in $cellsToClear) {
    $ws.Range($addr).ClearContents()
}

# ---------------------------------------------------------------------------
# 3. Write the corrected marker cells (rows 2-14).
# ---------------------------------------------------------------------------
$cellsToSet = @(
    @{ addr = "J2"; v = 2 },

    @{ addr = "I3"; v = 0 },
    @{ addr = "J3"; v = 2 },

    @{ addr = "G4"; v = 1 },
    @{ addr = "H4"; v = 2 },

    @{ addr = "D5"; v = 0 },
    @{ addr = "E5"; v = 1 },
    @{ addr = "F5"; v = 2 },

    @{ addr = "E6"; v = 0 },
    @{ addr = "F6"; v = 2 },

    @{ addr = "G7"; v = 0 },
    @{ addr = "H7"; v = 2 },

    @{ addr = "G8"; v = 0 },
    @{ addr = "H8"; v = 1 },
    @{ addr = "I8"; v = 2 },

    @{ addr = "F9"; v = 1 },
    @{ addr = "G9"; v = 2 },

    @{ addr = "G10"; v = 0 },
    @{ addr = "I10"; v = 1 },
    @{ addr = "J10"; v = 2 },

    @{ addr = "E11"; v = 0 },
    @{ addr = "F11"; v = 2 },

    @{ addr = "D12"; v = 0 },
    @{ addr = "E12"; v = 1 },
    @{ addr = "F12"; v = 2 },

    @{ addr = "G13"; v = 0 },
    @{ addr = "H13"; v = 2 },

    @{ addr = "G14"; v = 0 },
    @{ addr = "H14"; v = 1 },
    @{ addr = "I14"; v = 2 }
)
foreach ($item in $cellsToSet) {
    $ws.Range($item.addr).Value = $item.v
}

# ---------------------------------------------------------------------------
# 4. Drop the classes that no longer exist (old rows 15, 16, 17) -- this also
#    shrinks the sheet dimension from A1:P17 down to A1:P14.
# ---------------------------------------------------------------------------
$ws.Rows("15:17").Delete()

# ---------------------------------------------------------------------------
# 5. Restore the active-cell selection used when the fix was saved.
# ---------------------------------------------------------------------------
$ws.Range("J6").Select()
